$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.480300000000002
$ws.Range("A3").Value = -21.8559
$ws.Range("D3").Value = -7.375799999999995
$ws.Range("D12").Value = -7.302900000000001
$ws.Range("A14").Value = -21.897
$ws.Range("A21").Value = -20.22179999999998
$ws.Range("A23").Value = -20.08229999999998
$ws.Range("D24").Value = -7.3143
$ws.Range("A25").Value = -21.81089999999999
$ws.Range("B25").Value = 5.545600000000003
$ws.Range("D25").Value = -8.431899999999994
$ws.Range("A26").Value = -21.08769999999996
$ws.Range("B27").Value = 6.014700000000006
$ws.Range("A29").Value = -20.67059999999998
$ws.Range("B31").Value = 5.610300000000001
$ws.Range("B39").Value = 9.659200000000002
$ws.Range("B48").Value = 5.204800000000002
$ws.Range("D50").Value = -8.179200000000002
$ws.Range("B51").Value = 5.787400000000002
$ws.Range("B52").Value = 5.276700000000001
$ws.Range("A53").Value = -22.24240000000001
$ws.Range("D53").Value = -6.260299999999998
$ws.Range("B55").Value = 6.025399999999999
$ws.Range("B56").Value = 5.206099999999998
$ws.Range("A57").Value = -22.47620000000001
$ws.Range("B57").Value = 4.886499999999997
$ws.Range("D57").Value = -8.670400000000006
$ws.Range("A59").Value = -22.38650000000001
$ws.Range("D61").Value = -7.847299999999999
$ws.Range("D63").Value = -8.031299999999996
$ws.Range("A69").Value = -21.63049999999999
$ws.Range("D70").Value = -7.292999999999994
$ws.Range("B73").Value = 8.315399999999993
$ws.Range("A79").Value = -20.51260000000001
$ws.Range("A83").Value = -22.00290000000001
$ws.Range("D86").Value = -7.510199999999997
$ws.Range("B89").Value = 4.569599999999995
$ws.Range("B90").Value = 5.532700000000001
$ws.Range("A91").Value = -21.4835
$ws.Range("B92").Value = 4.758499999999994
$ws.Range("A93").Value = -20.74449999999998
$ws.Range("D98").Value = -8.367899999999999
$ws.Range("D100").Value = -8.133500000000003
$ws.Range("D102").Value = -7.596099999999996
